$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 26000
$ws.Range("E8").Value = 25100
$ws.Range("F8").Value = 27000
$ws.Range("G8").Value = 18100
$ws.Range("H8").Value = 13100
$ws.Range("I8").Value = 9800
$ws.Range("J8").Value = 8100

$ws.Range("D9").Value = 15500
$ws.Range("E9").Value = 17300
$ws.Range("F9").Value = 19200
$ws.Range("G9").Value = 11000
$ws.Range("H9").Value = 7900
$ws.Range("I9").Value = 5600
$ws.Range("J9").Value = 4900

$ws.Range("D10").Value = 10500
$ws.Range("E10").Value = 7700
$ws.Range("F10").Value = 7800
$ws.Range("G10").Value = 7100
$ws.Range("H10").Value = 5200
$ws.Range("I10").Value = 4200
$ws.Range("J10").Value = 3300

$ws.Range("D12").Value = 6200
$ws.Range("E12").Value = 6400
$ws.Range("F12").Value = 6100
$ws.Range("G12").Value = 4500
$ws.Range("H12").Value = 3000

$ws.Range("D17").Value = 35700
$ws.Range("E17").Value = 37500
$ws.Range("F17").Value = 37600
$ws.Range("G17").Value = 22800
$ws.Range("H17").Value = 15400
$ws.Range("I17").Value = 9000
$ws.Range("J17").Value = 7600

$ws.Range("D18").Value = -9700
$ws.Range("E18").Value = -12500
$ws.Range("F18").Value = -10600
$ws.Range("G18").Value = -4700
$ws.Range("H18").Value = -2300

$ws.Range("J20").Value = -400

$ws.Range("D21").Value = -5900
$ws.Range("E21").Value = -9800
$ws.Range("F21").Value = -7300
$ws.Range("G21").Value = -2500
$ws.Range("I21").Value = 1900
$ws.Range("J21").Value = 1500

$ws.Range("D23").Value = -9500
$ws.Range("E23").Value = -12700
$ws.Range("F23").Value = -10700
$ws.Range("G23").Value = -4900
$ws.Range("H23").Value = -2600

$ws.Range("D26").Value = -9600
$ws.Range("E26").Value = -12700
$ws.Range("F26").Value = -10800
$ws.Range("G26").Value = -4900
$ws.Range("H26").Value = -3000
$ws.Range("J26").Value = 0

$ws.Range("D27").Value = -9600
$ws.Range("E27").Value = -12700
$ws.Range("F27").Value = -10800
$ws.Range("G27").Value = -4900
$ws.Range("H27").Value = -3000
$ws.Range("J27").Value = 0

$ws.Range("J32").Value = 400

$ws.Range("D33").Value = -9600
$ws.Range("E33").Value = -12700
$ws.Range("F33").Value = -10800
$ws.Range("G33").Value = -4900
$ws.Range("H33").Value = -3000
$ws.Range("J33").Value = 0

$ws.Range("D35").Value = -9600
$ws.Range("E35").Value = -12700
$ws.Range("F35").Value = -10800
$ws.Range("G35").Value = -4900
$ws.Range("H35").Value = -3000
$ws.Range("J35").Value = 0

$ws.Range("D41").Value = 8500
$ws.Range("E41").Value = 8800
$ws.Range("F41").Value = 2300
$ws.Range("G41").Value = 9000
$ws.Range("H41").Value = 37500
$ws.Range("I41").Value = 300

$ws.Range("D42").Value = 15800
$ws.Range("E42").Value = 14100
$ws.Range("F42").Value = 35600
$ws.Range("G42").Value = 46200
$ws.Range("H42").Value = 800

$ws.Range("D43").Value = 5700
$ws.Range("E43").Value = 4600
$ws.Range("F43").Value = 3800
$ws.Range("G43").Value = 3600

$ws.Range("D44").Value = 10400
$ws.Range("E44").Value = 12600
$ws.Range("F44").Value = 8800
$ws.Range("G44").Value = 5900
$ws.Range("H44").Value = 4100
$ws.Range("I44").Value = 3100
$ws.Range("J44").Value = 2300

$ws.Range("D45").Value = 1700
$ws.Range("E45").Value = 1900
$ws.Range("H45").Value = 1100

$ws.Range("D46").Value = 42100
$ws.Range("E46").Value = 42100
$ws.Range("F46").Value = 52200
$ws.Range("G46").Value = 65600
$ws.Range("H46").Value = 44900
$ws.Range("I46").Value = 4600
$ws.Range("J46").Value = 3900

$ws.Range("D47").Value = 400

$ws.Range("D48").Value = 62400
$ws.Range("E48").Value = 26400
$ws.Range("F48").Value = 24000
$ws.Range("G48").Value = 21800
$ws.Range("H48").Value = 18300
$ws.Range("I48").Value = 5900
$ws.Range("J48").Value = 6000

$ws.Range("D49").Value = 2500
$ws.Range("E49").Value = 900
$ws.Range("F49").Value = 2100
$ws.Range("G49").Value = 3200

$ws.Range("D54").Value = 75200
$ws.Range("E54").Value = 69700
$ws.Range("F54").Value = 78700
$ws.Range("G54").Value = 91000
$ws.Range("H54").Value = 65000
$ws.Range("I54").Value = 12000
$ws.Range("J54").Value = 11000

$ws.Range("D57").Value = 3400
$ws.Range("E57").Value = 2000
$ws.Range("F57").Value = 2000
$ws.Range("G57").Value = 2600
$ws.Range("H57").Value = 1700
$ws.Range("I57").Value = 600

$ws.Range("D58").Value = 1300
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"

$ws.Range("D59").Value = 2600
$ws.Range("E59").Value = 4200
$ws.Range("F59").Value = 5200
$ws.Range("G59").Value = 3600
$ws.Range("H59").Value = 6300
$ws.Range("I59").Value = 4200
$ws.Range("J59").Value = 2900

$ws.Range("D60").Value = 7400
$ws.Range("E60").Value = 6200
$ws.Range("F60").Value = 7200
$ws.Range("G60").Value = 6200
$ws.Range("H60").Value = 8000
$ws.Range("I60").Value = 4800
$ws.Range("J60").Value = 3500

$ws.Range("D61").Value = 18400
$ws.Range("E61").Value = 5400
$ws.Range("F61").Value = 1400
$ws.Range("G61").Value = 2500
$ws.Range("H61").Value = 4300
$ws.Range("I61").Value = 4600
$ws.Range("J61").Value = 5200

$ws.Range("D62").Value = 100
$ws.Range("G62").Value = 2200
$ws.Range("I62").Value = 1200
$ws.Range("J62").Value = 1100

$ws.Range("D66").Value = 26000
$ws.Range("E66").Value = 12000
$ws.Range("F66").Value = 9700
$ws.Range("G66").Value = 11000
$ws.Range("H66").Value = 14000
$ws.Range("I66").Value = 10700
$ws.Range("J66").Value = 9800

$ws.Range("D72").Value = 43400
$ws.Range("E72").Value = 52600
$ws.Range("F72").Value = 65100
$ws.Range("G72").Value = 75800
$ws.Range("H72").Value = 47400
$ws.Range("I72").Value = 200

$ws.Range("D76").Value = 49200
$ws.Range("E76").Value = 57700
$ws.Range("F76").Value = 69000
$ws.Range("G76").Value = 80000
$ws.Range("H76").Value = 50900
$ws.Range("J76").Value = 1100

$ws.Range("D81").Value = -9600
$ws.Range("E81").Value = -12700
$ws.Range("F81").Value = -10800
$ws.Range("G81").Value = -4900
$ws.Range("H81").Value = -3000
$ws.Range("J81").Value = 0

$ws.Range("D83").Value = 3500
$ws.Range("E83").Value = 2900
$ws.Range("F83").Value = 3300
$ws.Range("G83").Value = 2400
$ws.Range("H83").Value = 1700
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 1400

$ws.Range("D89").Value = -7900
$ws.Range("E89").Value = -15000
$ws.Range("F89").Value = -13400
$ws.Range("G89").Value = -5600
$ws.Range("H89").Value = -1800

$ws.Range("D91").Value = -3800
$ws.Range("E91").Value = -4200
$ws.Range("G91").Value = -3000
$ws.Range("H91").Value = -12500
$ws.Range("J91").Value = -200

$ws.Range("D94").Value = -5300
$ws.Range("E94").Value = 17400
$ws.Range("F94").Value = 8000
$ws.Range("G94").Value = -52800
$ws.Range("H94").Value = -12800
$ws.Range("J94").Value = -300

$ws.Range("D100").Value = 12700
$ws.Range("E100").Value = 4100
$ws.Range("F100").Value = -1200
$ws.Range("G100").Value = 29800
$ws.Range("H100").Value = 51900

$ws.Range("E102").Value = 6500
$ws.Range("F102").Value = -6700
$ws.Range("G102").Value = -28500
$ws.Range("H102").Value = 37200
